$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.494.51"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.727.35"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'245.14"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4806"
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("D8").Value = "'0.2676"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").Value = "'0.06212"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "1.727.43"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").Value = "'0.07151"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "'15.68"
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("D13").Value = "'0.6187"
$ws.Range("E13").Value = "  +4.90%  "
$ws.Range("D14").Value = "'4.511"
$ws.Range("D15").Value = "'77.16"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "'0.9998"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "26.507.59"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'0.9997"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'0.000006928"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "'11.66"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").Value = "1.949.08"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'4.532"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "'8.955"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").Value = "'5.290"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "'136.50"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'15.32"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").Value = "'1.800"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("D28").Value = "'1.406"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'106.83"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'3.976"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").Value = "'0.08011"
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").Value = "'0.04560"
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.614"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6359"
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9896"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.9350"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.087"
$ws.Range("E38").Value = "  +9.47%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.419"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").Value = "'105.33"
$ws.Range("E40").Value = "  -7.66%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01502"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.677"
$ws.Range("E43").Value = "  +7.59%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3903"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.908"
$ws.Range("E45").Value = "  +10.72%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1190"
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05328"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'31.01"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.875"
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.268"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3432"
$ws.Range("E51").Value = "  +1.71%  "
